$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-07 03:07:01"

$wsMain = $wb.Worksheets.Item("Главные")
$wsLinear = $wb.Worksheets.Item("Линейные")

# ---- Update data cells on "Главные" (sheet2) ----

# Row 8: Gamaley Evgeniy
$wsMain.Range("C8").Value = 30
$wsMain.Range("D8").Value = 587
$wsMain.Range("E8").Value = 298
$wsMain.Range("F8").Value = 289
$wsMain.Range("G8").Value = 19.57
$wsMain.Range("H8").Value = 9.93
$wsMain.Range("I8").Value = 9.630000000000001
$wsMain.Range("J8").Value = 124
$wsMain.Range("K8").Value = 122
$wsMain.Range("L8").Value = 4
$wsMain.Range("M8").Value = 5
$wsMain.Range("N8").Value = 1
$wsMain.Range("Q8").Value = 1
$wsMain.Range("V8").Value = 16
$wsMain.Range("X8").Value = 2

# Row 11: Dudarov Aleksandr
$wsMain.Range("C11").Value = 25
$wsMain.Range("D11").Value = 582
$wsMain.Range("E11").Value = 270
$wsMain.Range("F11").Value = 312
$wsMain.Range("G11").Value = 23.28
$wsMain.Range("H11").Value = 10.8
$wsMain.Range("I11").Value = 12.48
$wsMain.Range("J11").Value = 120
$wsMain.Range("K11").Value = 111

# Row 18: Naumov Denis
$wsMain.Range("C18").Value = 31
$wsMain.Range("D18").Value = 541
$wsMain.Range("E18").Value = 260
$wsMain.Range("F18").Value = 281
$wsMain.Range("G18").Value = 17.45
$wsMain.Range("H18").Value = 8.390000000000001
$wsMain.Range("I18").Value = 9.06
$wsMain.Range("J18").Value = 100
$wsMain.Range("K18").Value = 118
$wsMain.Range("L18").Value = 4
$wsMain.Range("M18").Value = 3
$wsMain.Range("N18").Value = 2
$wsMain.Range("Q18").Value = 1
$wsMain.Range("V18").Value = 10
$wsMain.Range("X18").Value = 2

# Row 21: Romasko Evgeniy
$wsMain.Range("C21").Value = 29
$wsMain.Range("D21").Value = 412
$wsMain.Range("E21").Value = 186
$wsMain.Range("F21").Value = 226
$wsMain.Range("G21").Value = 14.21
$wsMain.Range("H21").Value = 6.41
$wsMain.Range("I21").Value = 7.79
$wsMain.Range("J21").Value = 83
$wsMain.Range("K21").Value = 98

# ---- Update data cells on "Линейные" (sheet3) ----

# Row 3: Bersenyov Maksim
$wsLinear.Range("C3").Value = 32
$wsLinear.Range("D3").Value = 534
$wsLinear.Range("E3").Value = 276
$wsLinear.Range("F3").Value = 258
$wsLinear.Range("G3").Value = 16.69
$wsLinear.Range("H3").Value = 8.630000000000001
$wsLinear.Range("I3").Value = 8.06
$wsLinear.Range("J3").Value = 118
$wsLinear.Range("K3").Value = 99
$wsLinear.Range("L3").Value = 2
$wsLinear.Range("M3").Value = 4
$wsLinear.Range("N3").Value = 1
$wsLinear.Range("Q3").Value = 2
$wsLinear.Range("V3").Value = 16
$wsLinear.Range("X3").Value = 2

# ---- Refresh as_of_utc timestamp (column AA) for every data row (2-26) on both sheets ----

for ($r = 2; $r -le 26; $r++) {
    $wsMain.Range("AA$r").Value = $newTimestamp
    $wsLinear.Range("AA$r").Value = $newTimestamp
}
